$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-12-18 23:55:24", -0.1208223765908042, -0.001971482140213996, 0.009527966303479201),
    @("2024-12-18 23:55:24", -0.121821148028539, -0.001779200443179996, 0.008669769622442898),
    @("2024-12-18 23:55:25", -0.1218764545902377, -0.001841303534061996, 0.008976461862238043)
)

$startRow = 98
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
